$wb = $excel.ActiveWorkbook

# --- Sheet bookkeeping: rename Sheet1 -> add, insert a new "sub" sheet after it ---
$wsAdd = $wb.Worksheets.Item("Sheet1")
$wsAdd.Name = "add"
$wsSub = $wb.Worksheets.Add($null, $wsAdd)
$wsSub.Name = "sub"

# --- "add" sheet: new op column with "+" operator, numeric A column ---
$wsAdd.Range("C1").Value = "op"
$wsAdd.Range("A2").Value = 1
$wsAdd.Range("C2").Value = "+"
$wsAdd.Range("A3").Value = 2
$wsAdd.Range("C3").Value = "+"

$wsAdd.PageSetup.PaperSize = 9
$wsAdd.PageSetup.Orientation = 1

$wsAdd.Range("A1:C3").Select() | Out-Null

# --- "sub" sheet: mirror layout with "-" operator data ---
$wsSub.Range("A1").Value = "a"
$wsSub.Range("B1").Value = "b"
$wsSub.Range("C1").Value = "op"
$wsSub.Range("A2").Value = 10
$wsSub.Range("B2").Value = 4
$wsSub.Range("C2").Value = "-"
$wsSub.Range("A3").Value = 12
$wsSub.Range("B3").Value = 5
$wsSub.Range("C3").Value = "-"

$wsSub.PageSetup.LeftMargin = 54
$wsSub.PageSetup.RightMargin = 54
$wsSub.PageSetup.TopMargin = 72
$wsSub.PageSetup.BottomMargin = 72
$wsSub.PageSetup.HeaderMargin = 36
$wsSub.PageSetup.FooterMargin = 36
$wsSub.PageSetup.PaperSize = 9
$wsSub.PageSetup.Orientation = 1

$wsSub.Range("B3").Select() | Out-Null
